$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of test-mail data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A45").Value = "Bestel je 100 M5-bouten zodra je kan?"
$logs.Range("B45").Value = "mailmind.test@zohomail.eu"
$logs.Range("C45").Value = "Testmail #4: Bestel je 100 M5-bouten zodra je kan?"
$logs.Range("D45").Value = "Inkoop / Bestellingen"
$logs.Range("E45").Value = "Geachte klant,`nDank voor uw e-mail. Het lijkt erop dat u ons per abuis heeft gecontacteerd. We willen u vriendelijk verzoeken ons te voorzien van wat meer informatie, zodat we u beter van dienst kunnen zijn. Kunt u ons meer vertellen over uw specifieke behoeften en het product waar u naar op zoek bent?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F45").Value = "2025-08-05 19:28:36"
$logs.Range("G45").Value = "Ja"
$logs.Range("H45").Value = "Nee"
$logs.Range("I45").Value = "Ja"
$logs.Range("J45").Value = "Nee"

# Multi-line content in E45 makes Excel flag a custom row height; AutoFit
# brings it back to the (un-flagged) default so row 45 serialises plainly,
# like every other row in the sheet.
$logs.Rows.Item(45).AutoFit()

# --- Extend the conditional-formatting ranges to include the new row ---
$logs.Range("D2:D44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D45"))
$logs.Range("G2:G44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G45"))
$logs.Range("H2:H44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H45"))
$logs.Range("I2:I44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I45"))
$logs.Range("J2:J44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J45"))

# --- Sheet "Dashboard": bump the "Inkoop / Bestellingen" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 6
